$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel; force text format so they stay as text
$textForceCells = @("D5", "D7", "D10", "D12", "D14", "D15", "D19", "D20", "D22", "D23", "D25", "D26", "D29", "D31", "D33", "D36", "D38", "D39", "D40", "D41", "D45", "D46", "D47", "D49")
foreach ($c in $textForceCells) { $ws.Range($c).NumberFormat = "@" }

# Apply new values
$ws.Range("D5").Value = "234.32"
$ws.Range("D7").Value = "58.63"
$ws.Range("D10").Value = "0.0786"
$ws.Range("D12").Value = "14.95"
$ws.Range("D14").Value = "21.07"
$ws.Range("D15").Value = "0.774"
$ws.Range("D19").Value = "6.17"
$ws.Range("D20").Value = "71.23"
$ws.Range("D22").Value = "228.59"
$ws.Range("D23").Value = "1.00"
$ws.Range("D25").Value = "2.40"
$ws.Range("D26").Value = "169.61"
$ws.Range("D29").Value = "19.52"
$ws.Range("D31").Value = "0.121"
$ws.Range("D33").Value = "0.0633"
$ws.Range("D36").Value = "1.83"
$ws.Range("D38").Value = "0.999"
$ws.Range("D39").Value = "5.40"
$ws.Range("D40").Value = "0.0977"
$ws.Range("D41").Value = "98.41"
$ws.Range("D45").Value = "4.31"
$ws.Range("D46").Value = "16.56"
$ws.Range("D47").Value = "1.16"
$ws.Range("D49").Value = "7.42"

$ws.Range("D2").Value = "37.768.06"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.084.29"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.391.71"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").Value = "2.080.73"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "37.705.10"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("E43").Value = "  -2.57%  "
$ws.Range("D44").Value = "1.457.89"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("E45").Value = "  +4.25%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "2.276.67"
$ws.Range("E51").Value = "  -1.26%  "

# Restore default style for the text-forced cells (remove explicit @ format)
foreach ($c in $textForceCells) { $ws.Range($c).Style = "Normal" }

